$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 32: Description (col B) is entered first, then Rule (col A),
# so that the shared-strings table receives them in that order.
$ws.Range("B32").Value = "Print the actual structure into the specified configuration file"
$ws.Range("A32").Value = "R `$ f `$ config.xml"

# Reuse the existing bordered look from row 30, then give A32 a brand new
# fill color (pink / FFFF66FF) distinct from the other rule highlights.
$ws.Range("A30").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A32").Interior.Color = 16738047

$ws.Range("B30").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection/view down to the newly added row, mirroring the
# author's navigation to the bottom of the sheet after the edit.
$ws.Range("A33").Select()
